# Update version 1.6.0 PREVIEW
# Renames the client test-data company names from "0407A" to "0507A"
# on the "Client" worksheet (column B: COMPANY_NAME, rows 2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Client")

$ws.Range("B2").Value = "Anh Tester Client 0507A1"
$ws.Range("B3").Value = "Anh Tester Client 0507A2"
$ws.Range("B4").Value = "Anh Tester Client 0507A3"
